# "Get Company Info" - populate extra company-detail columns (CEO, Media
# Contact, Auditor, Solicitor, CFO, Website) for the first few NZX
# companies, with a new header row, wrapped "Media Contact" column and
# a couple of widened columns to make the extra data readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "CompanyName"
$ws.Range("B1").Value = "NzxCompanyUrl"
$ws.Range("C1").Value = "CEO"
$ws.Range("D1").Value = "Media Contact"
$ws.Range("E1").Value = "Auditor"
$ws.Range("F1").Value = "Solicitor"
$ws.Range("G1").Value = "CFO"
$ws.Range("H1").Value = "Website"

# --- Ascension Capital Limited ---------------------------------------
$ws.Range("A2").Value = "Ascension Capital Limited"
$ws.Range("B2").Value = "/companies/ACE"
$ws.Range("D2").Value = "`nJohn Cilliers`nC/o Duncan Cotterill`nLevel 2, Tower Building`n50 Customhouse Quay`nWellington 6143`n+64 9 520 1020`nAscension Capital Limited website"
$ws.Range("E2").Value = "BDO Wellington"
$ws.Range("F2").Value = "Duncan Cotterill"

# --- AFC Group Holdings Limited ---------------------------------------
$ws.Range("A3").Value = "AFC Group Holdings Limited"
$ws.Range("B3").Value = "/companies/AFC"
$ws.Range("D3").Value = "`nPO Box 230122`nBotany`nAuckland`n+6499300245`nAFC Group Holdings Limited website"
$ws.Range("E3").Value = "William Buck"
$ws.Range("F3").Value = "DLA Piper"
$ws.Range("G3").Value = "Hao Long"
$ws.Range("H3").Value = "http://www.afcnz.com"

# --- Australian Foundation Investment Company Limited -----------------
$ws.Range("A4").Value = "Australian Foundation Investment Company Limited"
$ws.Range("B4").Value = "/companies/AFI"
$ws.Range("C4").Value = "Mark Freeman"
$ws.Range("D4").Value = "`nAndrew Porter`nMail Box 146,`n101 Collins Street,`nMelbourne, VIC 3000`n+6139650 9911`nAustralian Foundation Investment Company Limited website"
$ws.Range("E4").Value = "PriceWaterhouseCoopers"
$ws.Range("G4").Value = "Andrew Porter"
$ws.Range("H4").Value = "http://www.afi.com.au/"

# --- AFT Pharmaceuticals Limited ---------------------------------------
$ws.Range("A5").Value = "AFT Pharmaceuticals Limited"
$ws.Range("B5").Value = "/companies/AFT"
$ws.Range("C5").Value = "Hartley Atkinson"
$ws.Range("D5").Value = "`nMalcolm Tubby`nLevel 1`n129 Hurstmere Road`nTakapuna`nAuckland 0622`n+64 9 488 0232`nAFT Pharmaceuticals Limited website"
$ws.Range("E5").Value = "Deloitte"
$ws.Range("F5").Value = "Harmos Horton Lusk Limited"
$ws.Range("G5").Value = "Malcolm Tubby"
$ws.Range("H5").Value = "http://www.aftpharm.com"

# --- Smartshares Global Aggregate Bond ETF -----------------------------
$ws.Range("A6").Value = "Smartshares Global Aggregate Bond ETF"
$ws.Range("B6").Value = "/companies/AGG"
$ws.Range("D6").Value = "`nJohn McLean`nPO Box 2959 Wellington 6140`n0800808780`nSmartshares Global Aggregate Bond ETF website"
$ws.Range("E6").Value = "KPMG"
$ws.Range("F6").Value = "DLA Piper"
$ws.Range("H6").Value = "https://smartshares.co.nz/"

# --- Formatting ---------------------------------------------------------
# Wrap the long "Media Contact" addresses so they stay readable.
$ws.Range("D1:D6").WrapText = $true

# Widen the Solicitor / CFO columns to fit their contents.
$ws.Columns.Item(6).ColumnWidth = 25.5
$ws.Columns.Item(7).ColumnWidth = 26.666667

# Row heights to fit the wrapped, multi-line "Media Contact" entries.
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 315
$ws.Rows.Item(3).RowHeight = 180
$ws.Rows.Item(4).RowHeight = 330
$ws.Rows.Item(5).RowHeight = 270
$ws.Rows.Item(6).RowHeight = 240

# Select the whole sheet (mirrors the "select all" state the workbook was
# left in), with E9 as the active cell.
$ws.Range("E9").Select()
$ws.Cells.Select()
